$wb = $excel.ActiveWorkbook

# Work on the "Tipos de Datos" sheet (2nd sheet), which lists field names/types.
$ws = $wb.Worksheets.Item("Tipos de Datos")
$ws.Activate()

# IMSI and MSISDN rows previously had type "INTEGER"; update them to "STRING"
# so the now-unused "INTEGER" shared string is removed from the workbook.
$ws.Range("B4").Value = "STRING"
$ws.Range("B5").Value = "STRING"

# Update the saved selection on this sheet to match the new state.
$ws.Range("C6").Select()
